# Updating Voila app to use S5
# Append eight new "bayes_net" combo rows (18-25) to Sheet1, covering the
# new S5 early_summer / late_summer / winter / spring periods for both
# the hindcast and forecast model types.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$comment = "Added 10.03.2020"

$rows = @(
    @{ Row=18; A="hindcast"; B="bayes_net"; C="http://meteo.unican.es/tds5/dodsC/Copernicus/SYSTEM5_ecmwf_Seasonal_25Members_SFC.ncml"; D="early_summer"; E="1993:2016"; F="5,6,7" },
    @{ Row=19; A="hindcast"; B="bayes_net"; C="http://meteo.unican.es/tds5/dodsC/Copernicus/SYSTEM5_ecmwf_Seasonal_25Members_SFC.ncml"; D="late_summer";  E="1993:2016"; F="8,9,10" },
    @{ Row=20; A="hindcast"; B="bayes_net"; C="http://meteo.unican.es/tds5/dodsC/Copernicus/SYSTEM5_ecmwf_Seasonal_25Members_SFC.ncml"; D="winter";       E="1993:2016"; F="11,12,1" },
    @{ Row=21; A="hindcast"; B="bayes_net"; C="http://meteo.unican.es/tds5/dodsC/Copernicus/SYSTEM5_ecmwf_Seasonal_25Members_SFC.ncml"; D="spring";       E="1993:2016"; F="2,3,4" },
    @{ Row=22; A="forecast"; B="bayes_net"; C="http://meteo.unican.es/tds5/dodsC/Copernicus/SYSTEM5_ecmwf_forecast_Seasonal_51Members_SFC.ncml"; D="early_summer"; E="2017:2019"; F="5,6,7" },
    @{ Row=23; A="forecast"; B="bayes_net"; C="http://meteo.unican.es/tds5/dodsC/Copernicus/SYSTEM5_ecmwf_forecast_Seasonal_51Members_SFC.ncml"; D="late_summer";  E="2017:2019"; F="8,9,10" },
    @{ Row=24; A="forecast"; B="bayes_net"; C="http://meteo.unican.es/tds5/dodsC/Copernicus/SYSTEM5_ecmwf_forecast_Seasonal_51Members_SFC.ncml"; D="winter";       E="2017:2019"; F="11,12,1" },
    @{ Row=25; A="forecast"; B="bayes_net"; C="http://meteo.unican.es/tds5/dodsC/Copernicus/SYSTEM5_ecmwf_forecast_Seasonal_51Members_SFC.ncml"; D="spring";       E="2017:2019"; F="2,3,4" }
)

# Fill column-by-column (matches how the shared-string table ends up
# ordered: early_summer, 5,6,7 / 8,9,10 / Added 10.03.2020).
foreach ($r in $rows) { $ws.Range("A$($r.Row)").Value2 = $r.A }
foreach ($r in $rows) { $ws.Range("B$($r.Row)").Value2 = $r.B }
foreach ($r in $rows) { $ws.Range("C$($r.Row)").Value2 = $r.C }
foreach ($r in $rows) { $ws.Range("D$($r.Row)").Value2 = $r.D }
foreach ($r in $rows) {
    $ws.Range("E$($r.Row)").NumberFormat = "@"
    $ws.Range("E$($r.Row)").Value2 = $r.E
}
foreach ($r in $rows) {
    $ws.Range("F$($r.Row)").NumberFormat = "@"
    $ws.Range("F$($r.Row)").Value2 = $r.F
}
foreach ($r in $rows) { $ws.Range("G$($r.Row)").Value2 = 1 }
foreach ($r in $rows) { $ws.Range("H$($r.Row)").Value2 = $comment }

# Column D/H widths grow to fit the newly added "early_summer" label and
# the shorter "Added 10.03.2020" comment text.
$ws.Columns("D").ColumnWidth = 11.17
$ws.Columns("H").ColumnWidth = 39

# Mirror the author's final cursor position/selection in the sheet.
$ws.Range("H28").Select() | Out-Null
